$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 value
$ws.Range("C2").Value = 23602299

# Update row 3 (ocap -> non_pdi)
$ws.Range("A3").Value = "non_pdi (5-17 y.o.)"
$ws.Range("B3").Value = "non_pdi"
$ws.Range("C3").Value = 12525253

# Update row 4 (idp -> pdi)
$ws.Range("A4").Value = "pdi (5-17 y.o.)"
$ws.Range("B4").Value = "pdi"
$ws.Range("C4").Value = 11077047

# Update row 5 (ret -> Girls)
$ws.Range("A5").Value = "Girls (5-17 y.o.)"
$ws.Range("B5").Value = "All population groups"
$ws.Range("C5").Value = 10381716

# Update row 6 (ndsp -> Boys)
$ws.Range("A6").Value = "Boys (5-17 y.o.)"
$ws.Range("B6").Value = "All population groups"
$ws.Range("C6").Value = 13220583

# Update row 7 (Girls -> ECE)
$ws.Range("A7").Value = "ECE (5 y.o.)"
$ws.Range("B7").Value = "All population groups"
$ws.Range("C7").Value = 2401923

# Update row 8 (Boys -> Children with disability)
$ws.Range("A8").Value = "Children with disability"
$ws.Range("B8").Value = "All population groups"
$ws.Range("C8").Value = 2360230

# Delete old rows 9 and 10 (ECE / Children with disability moved up into rows 7-8)
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()
